$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-02 Monday", "2024-09-03 Tuesday"),
    @("30-26=4", "75-30=45"),
    @("17+67=84", "76-13=63"),
    @("14+77=91", "78-11=67"),
    @("93-79=14", "75-70=5"),
    @("33-9=24", "37+28=65"),
    @("95-24=71", "31-28=3"),
    @("45-6=39", "56+38=94"),
    @("58+19=77", "30+36=66"),
    @("28+54=82", "37+49=86"),
    @("58-19=39", "82-72=10"),
    @("20+68=88", "45-27=18"),
    @("47+47=94", "73-56=17"),
    @("35-1=34", "12+57=69"),
    @("87-42=45", "18+30=48"),
    @("94+5=99", "61-5=56"),
    @("61-7=54", "43-2=41"),
    @("33+38=71", "89-72=17"),
    @("55+25=80", "73-67=6"),
    @("58-52=6", "0+83=83"),
    @("28+32=60", "8+70=78"),
    @("69+20=89", "11-2=9"),
    @("18+3=21", "52-31=21"),
    @("70+19=89", "18+66=84"),
    @("65-38=27", "15+5=20"),
    @("92-38=54", "73-3=70"),
    @("41+39=80", "55-33=22"),
    @("31+38=69", "88-51=37"),
    @("85-74=11", "91-17=74"),
    @("15-6=9", "13+10=23"),
    @("99-70=29", "87+2=89"),
    @("77-9=68", "49+5=54"),
    @("49+22=71", "30-3=27"),
    @("29-27=2", "17+40=57"),
    @("37+20=57", "50-12=38"),
    @("69-29=40", "85-72=13"),
    @("83-60=23", "10+0=10"),
    @("34+57=91", "48-27=21"),
    @("7+12=19", "45+52=97"),
    @("48-16=32", "65-8=57"),
    @("30-8=22", "90-50=40"),
    @("27+26=53", "21+70=91"),
    @("27+10=37", "27+62=89"),
    @("6+83=89", "49+27=76"),
    @("63+28=91", "57-51=6"),
    @("47+43=90", "57+12=69"),
    @("34-6=28", "25+31=56"),
    @("74-17=57", "45+17=62"),
    @("38+25=63", "54+45=99"),
    @("80-51=29", "22+7=29"),
    @("89-61=28", "82-53=29"),
    @("94-25=69", "55+29=84"),
    @("0+43=43", "85-39=46"),
    @("56+29=85", "87-80=7"),
    @("57-42=15", "82-12=70"),
    @("40+34=74", "4+1=5"),
    @("79+7=86", "22-14=8"),
    @("6-2=4", "77-12=65"),
    @("51+44=95", "55+7=62"),
    @("36+25=61", "79-35=44"),
    @("4+27=31", "37+53=90"),
    @("48+18=66", "38+57=95"),
    @("51-35=16", "47+1=48"),
    @("64-38=26", "8+39=47"),
    @("91+1=92", "80-17=63"),
    @("94-28=66", "39-24=15"),
    @("68-0=68", "46-45=1"),
    @("84-9=75", "0+82=82"),
    @("91-44=47", "68-17=51"),
    @("61-57=4", "58-30=28"),
    @("60-2=58", "84-51=33"),
    @("56+23=79", "10+31=41"),
    @("82+7=89", "22+19=41"),
    @("28+44=72", "37-16=21"),
    @("18+60=78", "16+15=31"),
    @("56+30=86", "39-36=3"),
    @("77-69=8", "43-2=41"),
    @("56-48=8", "7+18=25"),
    @("90-79=11", "38+45=83"),
    @("75-54=21", "20-3=17"),
    @("14+64=78", "54-33=21"),
    @("79-66=13", "66-20=46"),
    @("86-7=79", "69+5=74"),
    @("1+10=11", "61+34=95"),
    @("69+15=84", "32-21=11"),
    @("50+22=72", "1+25=26"),
    @("27+9=36", "95-77=18"),
    @("34+29=63", "42+33=75"),
    @("26+4=30", "28+61=89"),
    @("24-9=15", "81+12=93"),
    @("77+11=88", "42-40=2"),
    @("5+77=82", "5+65=70"),
    @("54+44=98", "81-25=56"),
    @("62+27=89", "95-73=22"),
    @("50-18=32", "54-45=9"),
    @("60+18=78", "64-57=7"),
    @("11+72=83", "18+15=33"),
    @("3+88=91", "82-53=29"),
    @("19+54=73", "12-12=0"),
    @("64-37=27", "68+0=68"),
    @("96-60=36", "23+55=78")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $find = $range.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
